# Create the "Account Funding Withdraw" form sheet, mirroring the layout
# of the existing "Account Funding Deposit" sheet, and make it the active
# (last, selected) tab of the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Add()
$ws.Name = "Account Funding Withdraw"

# Fill column A first, then column B (top-to-bottom), matching how the
# original author would have typed the form in so shared-string entries
# land in the same order.
$ws.Range("A1").Value = "public address"
$ws.Range("A2").Value = "0x5befc48f793f5f2595ca460f72ef785fe0f7c842"
$ws.Range("B1").Value = "amount to withdraw"
$ws.Range("B2").Value = 1

# ColumnWidth is stored internally on a whole-pixel grid (padding + 6px per
# character), so these inputs are chosen to land the persisted <col width>
# as close as the grid allows to the source sheet's 79.7109375 / 25.85546875.
$ws.Columns.Item(1).ColumnWidth = 78.8333333333333
$ws.Columns.Item(2).ColumnWidth = 25

# Move the new sheet to the end of the tab strip (Worksheets.Add() inserts
# it before the workbook's current active sheet, so it needs relocating).
$ws.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-resolve the sheet by name after the move and make it the active /
# selected tab, with the same cell selection the source file shipped with.
$ws = $wb.Worksheets.Item("Account Funding Withdraw")
$ws.Activate()
$ws.Range("C8").Select()
